# Applies per-row Price (D) and Volume(1h) (E) updates for the cryptos list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.617.07'
$ws.Range("E2").Value = '  -1.77%  '
$ws.Range("D3").Value = '2.905.45'
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("E4").Value = '  +0.00%  '
$r = $ws.Range("D5")
$r.Formula = '="529.10"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E5").Value = '  -2.07%  '
$r = $ws.Range("D6")
$r.Formula = '="143.98"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E6").Value = '  -5.74%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -1.10%  '
$ws.Range("D9").Value = '2.913.88'
$ws.Range("E9").Value = '  -2.01%  '
$ws.Range("E10").Value = '  -3.30%  '
$ws.Range("E11").Value = '  -0.76%  '
$r = $ws.Range("D12")
$r.Formula = '="0.362"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("D13").Value = '3.415.10'
$ws.Range("E13").Value = '  -1.92%  '
$ws.Range("E14").Value = '  +2.57%  '
$ws.Range("D15").Value = '60.611.59'
$ws.Range("E15").Value = '  -1.82%  '
$r = $ws.Range("D16")
$r.Formula = '="22.78"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E16").Value = '  -3.69%  '
$ws.Range("D17").Value = '2.908.88'
$ws.Range("E17").Value = '  -2.16%  '
$ws.Range("E18").Value = '  -3.43%  '
$r = $ws.Range("D19")
$r.Formula = '="5.05"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E19").Value = '  -1.35%  '
$r = $ws.Range("D20")
$r.Formula = '="11.70"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E20").Value = '  -2.04%  '
$r = $ws.Range("D21")
$r.Formula = '="362.24"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E21").Value = '  -4.65%  '
$r = $ws.Range("D22")
$r.Formula = '="6.63"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("E23").Value = '  -0.04%  '
$r = $ws.Range("D24")
$r.Formula = '="5.68"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E24").Value = '  -0.08%  '
$r = $ws.Range("D25")
$r.Formula = '="64.83"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E25").Value = '  -0.34%  '
$r = $ws.Range("D26")
$r.Formula = '="0.456"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E26").Value = '  -3.05%  '
$ws.Range("E27").Value = '  -2.95%  '
$ws.Range("E28").Value = '  +0.15%  '
$r = $ws.Range("D29")
$r.Formula = '="7.86"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E29").Value = '  -4.52%  '
$ws.Range("E30").Value = '  -6.98%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("E32").Value = '  -2.06%  '
$r = $ws.Range("D33")
$r.Formula = '="19.78"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E33").Value = '  -3.18%  '
$r = $ws.Range("D34")
$r.Formula = '="152.30"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E34").Value = '  -4.08%  '
$r = $ws.Range("D35")
$r.Formula = '="4.39"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E35").Value = '  -5.53%  '
$ws.Range("E36").Value = '  -5.88%  '
$ws.Range("E37").Value = '  -4.42%  '
$ws.Range("E38").Value = '  -5.34%  '
$r = $ws.Range("D39")
$r.Formula = '="37.69"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E39").Value = '  +1.59%  '
$ws.Range("E40").Value = '  -4.34%  '
$r = $ws.Range("D41")
$r.Formula = '="3.73"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E41").Value = '  -4.46%  '
$ws.Range("D42").Value = '2.298.22'
$ws.Range("E42").Value = '  -4.78%  '
$r = $ws.Range("D43")
$r.Formula = '="0.648"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E43").Value = '  -2.22%  '
$ws.Range("E44").Value = '  -1.27%  '
$r = $ws.Range("D45")
$r.Formula = '="20.58"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E45").Value = '  -7.39%  '
$ws.Range("E46").Value = '  +0.02%  '
$r = $ws.Range("D47")
$r.Formula = '="5.01"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E47").Value = '  +0.95%  '
$ws.Range("E48").Value = '  -2.69%  '
$r = $ws.Range("D49")
$r.Formula = '="10.32"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E49").Value = '  -1.47%  '
$r = $ws.Range("D50")
$r.Formula = '="0.0925"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E50").Value = '  -2.45%  '
$r = $ws.Range("D51")
$r.Formula = '="251.97"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range("E51").Value = '  -5.48%  '

$excel.CutCopyMode = 0

